$wb = $excel.ActiveWorkbook

# Updated counts for "F" column ("想去人数") on sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8087
$ws1.Range("F13").Value = 3318
$ws1.Range("F14").Value = 221
$ws1.Range("F15").Value = 109
$ws1.Range("F16").Value = 753
$ws1.Range("F21").Value = 281
$ws1.Range("F22").Value = 736
$ws1.Range("F23").Value = 353
$ws1.Range("F26").Value = 130
$ws1.Range("F32").Value = 582
$ws1.Range("F34").Value = 38

# Same updated counts mirrored on sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 8087
$ws4.Range("F16").Value = 3318
$ws4.Range("F17").Value = 221
$ws4.Range("F18").Value = 109
$ws4.Range("F20").Value = 753
$ws4.Range("F26").Value = 281
$ws4.Range("F27").Value = 736
$ws4.Range("F28").Value = 353
$ws4.Range("F31").Value = 130
$ws4.Range("F37").Value = 582
$ws4.Range("F39").Value = 38
